# Regenerate the s_vals data (filter save games) by updating the
# numeric stat columns (B:E) and the derived "sum" column (G) for
# rows 2-7. Column F (Win) is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    3 = @(3.286832544864788, 1.655778082260271, 261.3203778131603, 10.19245300693656,  276.4554414472219)
    4 = @(1.455362044514542, 1.655778082260271, 22.3905356188092,  0.4942365360607697, 25.99591228164478)
    5 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641,0.4942365360607697, 4.358119930609447)
    6 = @(0.6606524410359556,0.306821227259698, 0.7527432677738641,10.19245300693656,  11.91266994300607)
    7 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641,0.4942365360607697, 6.189590430959694)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - sum
}
